$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the full data range as Text so that numeric-looking values
# (gas-used counts and ETH amounts) are stored as text, as in the source.
$ws.Range("A1:F58").NumberFormat = "@"

$data = @(
    @('Contract','Name','Tx','Description','GasUsed','GasUsed (ETH) '),
    @('L2Registry','addManager','','','50991','0.00050991'),
    @('OperatorFactory','setAddresses','','','115592','0.00115592'),
    @('Layer2Manager','setAddresses','','','217793','0.00217793'),
    @('Layer2Manager','setMinimumInitialDepositAmount','','','57293','0.00057293'),
    @('SystemConfig','setAddresses','','','146738','0.00146738'),
    @('L2Registry','registerSystemConfigByManager','','','89648','0.00089648'),
    @('L2Registry','setAddresses','','','146750','0.0014675'),
    @('L2Registry','registerSystemConfigByManager','','Thanos SystemConfig ','98997','0.00098997'),
    @('Layer2Manager','registerLayer2Candidate','','using registerLayer2Candidate function','4862471','0.04862471'),
    @('SeigManagerV1_3','setLayer2StartBlock','','','56913','0.00056913'),
    @('DepositManager','approveAndCall','','deposit TON to Titan Candidate','374814','0.00374814'),
    @('DepositManager','deposit(address,uint256)','','deposit WTON to Titan Candidate','303991','0.00303991'),
    @('DepositManager','deposit(address,address,uint256)','','deposit WTON to Titan Candidate','254169','0.00254169'),
    @('SeigManagerV1_3','updateSeigniorageLayer(layerAddress)','','the first updateSeigniorage : no give seigniorage to l2','429640','0.0042964'),
    @('SeigManagerV1_3','updateSeigniorageLayer(layerAddress)','','the second updateSeigniorage : not operator','437487','0.00437487'),
    @('SeigManagerV1_3','updateSeigniorage(1)','','the third updateSeigniorage of operator with claiming','480766','0.00480766'),
    @('SeigManagerV1_3','updateSeigniorage(2)','','the forth updateSeigniorage of operator with staking ','588263','0.00588263'),
    @('DepositManager','requestWithdrawal','','','341275','0.00341275'),
    @('DepositManager','processRequest','','','197362','0.00197362'),
    @('DepositManager','approveAndCall','','deposit TON to DAOCandidate','315110','0.0031511'),
    @('DepositManager','approveAndCall','','deposit TON to DAOCandidate','356170','0.0035617'),
    @('DepositManager','deposit(address,uint256)','','deposit WTON to DAOCandidate','249099','0.00249099'),
    @('DepositManager','deposit(address,address,uint256)','','deposit WTON to DAOCandidate','290666','0.00290666'),
    @('SeigManagerV1_3','updateSeigniorageLayer(layerAddress)','','updateSeigniorage to DAOCandidate','347933','0.00347933'),
    @('DepositManager','requestWithdrawal(address,uint256)','','','278176','0.00278176'),
    @('DepositManager','processRequests(address,uint256,bool)','','','145889','0.00145889'),
    @('DepositManager','approveAndCall','','deposit TON to Titan Layer2 Candidate','357714','0.00357714'),
    @('DepositManager','approveAndCall','','deposit TON to DAO Candidate','322941','0.00322941'),
    @('DepositManager','withdrawAndDepositL2','','','476145','0.00476145'),
    @('L2Registry','rejectLayer2Candidate','','legacySystemConfig','115789','0.00115789'),
    @('L2Registry','registerLayer2Candidate','','thanos SystemConfig','4823637','0.04823637'),
    @('SeigManager','updateSeigniorageLayer','','titanLayerAddress','355902','0.00355902'),
    @('titanLayerContract','updateSeigniorage','','with claim','355815','0.00355815'),
    @('titanLayerContract','updateSeigniorage','','with staking','355815','0.00355815'),
    @('DepositManager','deposit(address,uint256)','','','286891','0.00286891'),
    @('SeigManager','updateSeigniorageLayer','','no give seigniorage to l2','435106','0.00435106'),
    @('DepositManager','deposit(address,address,uint256)','','','254169','0.00254169'),
    @('SeigManager','updateSeigniorageLayer','','give seigniorage to l2','408491','0.00408491'),
    @('SeigManager','updateSeigniorageLayer','','not operator','408491','0.00408491'),
    @('Layer2Contract','updateSeigniorage','','operator','485970','0.0048597'),
    @('Layer2Contract','updateSeigniorage','','with operator''s staking','593467','0.00593467'),
    @('DepositManager','requestWithdrawal','','','341275','0.00341275'),
    @('DepositManager','processRequest','','','180262','0.00180262'),
    @('L2Registry','restoreLayer2Candidate','','','72294','0.00072294'),
    @('SeigManager','updateSeigniorageLayer','','','425987','0.00425987'),
    @('SeigManager','updateSeigniorageLayer','','','403287','0.00403287'),
    @('SeigManager','updateSeigniorageLayer','','','403287','0.00403287'),
    @('Layer2Contract','updateSeigniorage','','with operator''s claim','463666','0.00463666'),
    @('Layer2Contract','updateSeigniorage','','with operator''s staking','567880','0.0056788'),
    @('DepositManager','requestWithdrawal','','','323583','0.00323583'),
    @('DepositManager','processRequest','','','128962','0.00128962'),
    @('TonContract','approveAndCall','','DepositManager.onApprove','322929','0.00322929'),
    @('DepositManager','deposit(address,uint256)','','','249099','0.00249099'),
    @('DepositManager','deposit(address,address,uint256)','','','257425','0.00257425'),
    @('SeigManager','updateSeigniorageLayer','','','354699','0.00354699'),
    @('DepositManager','requestWithdrawal','','','326832','0.00326832'),
    @('DepositManager','processRequests','','','109935','0.00109935')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

